$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add retrospective answers in row 2 (Start / Stop / Continue columns)
$ws.Range("A2").Value = "Starting development sooner"
$ws.Range("B2").Value = "N/A"
$ws.Range("C2").Value = "Working towards MVP (Minimal Viable Product)"

# Extra note in row 3
$ws.Range("A3").Value = "Potentially do UI Mockups"

# Row 2 holds two-line wrapped text, so it is taller than the default row
$ws.Rows.Item(2).RowHeight = 30

# Leave the same selection state the source workbook was saved with
$ws.Range("B3").Select() | Out-Null
